$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# Locate the "Arrow: Right 14" shape (PowerPoint shows it as "Right Arrow 14")
# and bring it to the front of the z-order, so it is drawn after (on top of /
# later in the XML than) the block of shapes that currently follows it.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 15) {
        $target = $sh
        break
    }
}

if ($target -ne $null) {
    $target.ZOrder(0)  # msoBringToFront
}
